$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.689.79"
$ws.Range("D3").Value = "1.584.62"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  +1.46%  "
$ws.Range("D5").Value = "'205.83"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("D8").Value = "'22.22"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.810.07"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "1.577.78"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").Value = "27.670.71"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "'219.00"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("D20").Value = "'7.30"
$ws.Range("E20").Value = "  -4.31%  "
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  -4.75%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'155.36"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'6.81"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "'15.11"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("D33").Value = "1.383.21"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'0.978"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "'0.537"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").Value = "'0.823"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Value = "'0.978"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").Value = "'63.48"
$ws.Range("D46").Value = "'5.23"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").Value = "1.720.50"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "'88.31"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +12.81%  "
$ws.Range("D50").Value = "'0.0973"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'0.0500"
$ws.Range("E51").Value = "  -0.50%  "
